$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 72 (weekly update adds the
# latest price record at the top of this product's block); this shifts
# the former rows 72-75 down to 73-76, preserving all their data.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly record. It
# shares the same product/market metadata as the rest of the block,
# only the date (D) and the price (J) are new for this week.
$ws.Range("A72").Value = 10
$ws.Range("B72").Value = "Vega Modelo de Temuco"
$ws.Range("C72").Value = "La Araucanía"
$ws.Range("D72").Value = 45021
$ws.Range("E72").Value = 9
$ws.Range("F72").Value = 100112010
$ws.Range("G72").Value = "Achicoria"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 25
$ws.Range("K72").Value = 10000
$ws.Range("L72").Value = 10000
$ws.Range("M72").Value = 10000
$ws.Range("N72").Value = "$/caja 18 unidades"
$ws.Range("O72").Value = "Región Metropolitana"
$ws.Range("P72").Value = 556
$ws.Range("Q72").Value = 18
$ws.Range("R72").Value = "Hortaliza"

# Carry over the date format from the rest of the column (style index 2,
# the date/time number format) onto the new cell.
$ws.Range("D72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
